# New weekly price observation: insert a row right after the header/latest
# block (new physical row 9) and push the previous rows 9-31 down to 10-32.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 9..31 down by one row, preserving formatting
# (the new blank row inherits the date style from the row above it).
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new observation.
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9, 3).Value = "Ñuble"
$ws.Cells.Item(9, 4).Value = 44972
$ws.Cells.Item(9, 5).Value = 16
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100101
$ws.Cells.Item(9, 8).Value = "Berries"
$ws.Cells.Item(9, 9).Value = 100101001
$ws.Cells.Item(9, 10).Value = "Arándano (blue)"
$ws.Cells.Item(9, 11).Value = "Sin especificar"
$ws.Cells.Item(9, 12).Value = "Segunda"
$ws.Cells.Item(9, 13).Value = 30
$ws.Cells.Item(9, 14).Value = 2500
$ws.Cells.Item(9, 15).Value = 2500
$ws.Cells.Item(9, 16).Value = 2500
$ws.Cells.Item(9, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(9, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(9, 19).Value = 1250
$ws.Cells.Item(9, 20).Value = 2
